$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "2FileTemplate"
$ws.Range("B8").Value = "\CONFIG\Template_kyocera_2.xlsx"

$ws.Range("B21").Select()
